# Fruta / hortaliza, semanal
# Weekly refresh: insert a new data row at row 55 (pushing every
# subsequent record down by one), so the last previously-existing
# record (old row 113) becomes the new row 114.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 55, shifting rows
# 55..113 down to 56..114.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A55").Value = 10
$ws.Range("B55").Value = "Vega Modelo de Temuco"
$ws.Range("C55").Value = "La Araucanía"
$ws.Range("D55").Value = 44902
$ws.Range("E55").Value = 9
$ws.Range("F55").Value = "Fruta"
$ws.Range("G55").Value = 100101
$ws.Range("H55").Value = "Berries"
$ws.Range("I55").Value = 100101001
$ws.Range("J55").Value = "Arándano (blue)"
$ws.Range("K55").Value = "Sin especificar"
$ws.Range("L55").Value = "Primera"
$ws.Range("M55").Value = 400
$ws.Range("N55").Value = 2200
$ws.Range("O55").Value = 2200
$ws.Range("P55").Value = 2200
$ws.Range("Q55").Value = "$/kilo"
$ws.Range("R55").Value = "Región del Maule"
$ws.Range("S55").Value = 2200
$ws.Range("T55").Value = 1
